# Insert a new weekly price record for "Macroferia Regional de Talca" (Acelga)
# right before the existing row 194. This pushes the old rows 194-304 down by
# one row (to 195-305), exactly mirroring an Excel "Insert Row" + fill-in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 194 - everything currently at/after
# row 194 shifts down by one (old row194 -> 195, ..., old row304 -> 305).
$ws.Rows.Item(194).Insert()

# Populate the newly inserted row 194 with the new record.
$ws.Range("A194").Value = 5
$ws.Range("B194").Value = "Macroferia Regional de Talca"
$ws.Range("C194").Value = "Maule"
$ws.Range("D194").Value = 44806
$ws.Range("E194").Value = 7
$ws.Range("F194").Value = 100112009
$ws.Range("G194").Value = "Acelga"
$ws.Range("H194").Value = "Sin especificar"
$ws.Range("I194").Value = "Primera"
$ws.Range("J194").Value = 500
$ws.Range("K194").Value = 2500
$ws.Range("L194").Value = 2500
$ws.Range("M194").Value = 2500
$ws.Range("N194").Value = '$/docena de atados (4 kilos)'
$ws.Range("O194").Value = "Región del Maule"
$ws.Range("P194").Value = 625
$ws.Range("Q194").Value = 4
$ws.Range("R194").Value = "Hortaliza"
